# Daily attendance processing - 2025-10-10 09:22:12
# Normalizes the "Recorded By" (column G) cell values: whenever the
# comma-separated list of recorders in a cell contains the literal
# "System" entry, the list order is reversed (e.g. "System, foo@bar.com"
# becomes "foo@bar.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ",\s*"

    if ($parts.Count -gt 1 -and ($parts -contains "System")) {
        $reversed = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $cell.Value = [string]::Join(", ", $reversed)
    }
}
